$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.686499
$ws.Range("H2").Value = 14.059497
$ws.Range("I2").Value = 0.05495643536242917
$ws.Range("J2").Value = 0.05495643536242918
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 33.42823600000001
$ws.Range("N2").Value = 100.284708
$ws.Range("O2").Value = 0.780497154124833
$ws.Range("P2").Value = 0.780497154124833
$ws.Range("Q2").Value = 156.661394585764
$ws.Range("R2").Value = 1409.952551271876
$ws.Range("S2").Value = 0.04289334140122131
$ws.Range("T2").Value = 0.04289334140122131

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.686499
$ws.Range("H3").Value = 14.059497
$ws.Range("I3").Value = 0.05495643536242917
$ws.Range("J3").Value = 0.05495643536242918
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.598076333333333
$ws.Range("N3").Value = 4.794229
$ws.Range("O3").Value = 0.0373125889813903
$ws.Range("P3").Value = 0.0373125889813903
$ws.Range("Q3").Value = 7.489383138090333
$ws.Range("R3").Value = 67.40444824281299
$ws.Range("S3").Value = 0.002050566884560663
$ws.Range("T3").Value = 0.002050566884560663

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.686499
$ws.Range("H4").Value = 14.059497
$ws.Range("I4").Value = 0.05495643536242917
$ws.Range("J4").Value = 0.05495643536242918
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.803102
$ws.Range("N4").Value = 23.409306
$ws.Range("O4").Value = 0.1821902568937766
$ws.Range("P4").Value = 0.1821902568937767
$ws.Range("Q4").Value = 36.569229719898
$ws.Range("R4").Value = 329.123067479082
$ws.Range("S4").Value = 0.0100125270766472
$ws.Range("T4").Value = 0.0100125270766472

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 46.42982133333334
$ws.Range("H5").Value = 139.289464
$ws.Range("I5").Value = 0.5444613292341401
$ws.Range("J5").Value = 0.5444613292341403
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.42823600000001
$ws.Range("N5").Value = 100.284708
$ws.Range("O5").Value = 0.780497154124833
$ws.Range("P5").Value = 0.780497154124833
$ws.Range("Q5").Value = 1552.067024968502
$ws.Range("R5").Value = 13968.60322471651
$ws.Range("S5").Value = 0.4249505179982701
$ws.Range("T5").Value = 0.4249505179982702

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 46.42982133333334
$ws.Range("H6").Value = 139.289464
$ws.Range("I6").Value = 0.5444613292341401
$ws.Range("J6").Value = 0.5444613292341403
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.598076333333333
$ws.Range("N6").Value = 4.794229
$ws.Range("O6").Value = 0.0373125889813903
$ws.Range("P6").Value = 0.0373125889813903
$ws.Range("Q6").Value = 74.1983986336951
$ws.Range("R6").Value = 667.7855877032559
$ws.Range("S6").Value = 0.02031526179397489
$ws.Range("T6").Value = 0.0203152617939749

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 46.42982133333334
$ws.Range("H7").Value = 139.289464
$ws.Range("I7").Value = 0.5444613292341401
$ws.Range("J7").Value = 0.5444613292341403
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.803102
$ws.Range("N7").Value = 23.409306
$ws.Range("O7").Value = 0.1821902568937766
$ws.Range("P7").Value = 0.1821902568937767
$ws.Range("Q7").Value = 362.296631705776
$ws.Range("R7").Value = 3260.669685351984
$ws.Range("S7").Value = 0.09919554944189508
$ws.Range("T7").Value = 0.09919554944189513

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 34.16029866666667
$ws.Range("H8").Value = 102.480896
$ws.Range("I8").Value = 0.4005822354034306
$ws.Range("J8").Value = 0.4005822354034306
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 33.42823600000001
$ws.Range("N8").Value = 100.284708
$ws.Range("O8").Value = 0.780497154124833
$ws.Range("P8").Value = 0.780497154124833
$ws.Range("Q8").Value = 1141.918525659819
$ws.Range("R8").Value = 10277.26673093837
$ws.Range("S8").Value = 0.3126532947253415
$ws.Range("T8").Value = 0.3126532947253415

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 34.16029866666667
$ws.Range("H9").Value = 102.480896
$ws.Range("I9").Value = 0.4005822354034306
$ws.Range("J9").Value = 0.4005822354034306
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.598076333333333
$ws.Range("N9").Value = 4.794229
$ws.Range("O9").Value = 0.0373125889813903
$ws.Range("P9").Value = 0.0373125889813903
$ws.Range("Q9").Value = 54.59076483879822
$ws.Range("R9").Value = 491.316883549184
$ws.Range("S9").Value = 0.01494676030285474
$ws.Range("T9").Value = 0.01494676030285474

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 34.16029866666667
$ws.Range("H10").Value = 102.480896
$ws.Range("I10").Value = 0.4005822354034306
$ws.Range("J10").Value = 0.4005822354034306
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.803102
$ws.Range("N10").Value = 23.409306
$ws.Range("O10").Value = 0.1821902568937766
$ws.Range("P10").Value = 0.1821902568937767
$ws.Range("Q10").Value = 266.556294846464
$ws.Range("R10").Value = 2399.006653618176
$ws.Range("S10").Value = 0.07298218037523434
$ws.Range("T10").Value = 0.07298218037523434
